$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2..302).
# Update all of them from 45190 (2023-09-21) to 45192 (2023-09-23).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 302 }

$ws.Range("C2:C$lastRow").Value = 45192
